$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "29.131.57"
$ws.Range("E2").Value2 = "  +2.53%  "
$ws.Range("D3").Value2 = "1.910.60"
$ws.Range("E3").Value2 = "  +2.12%  "
$ws.Range("D4").Value2 = "1.001"
$ws.Range("E4").Value2 = "  -0.03%  "
$ws.Range("D5").Value2 = "334.47"
$ws.Range("E5").Value2 = "  -1.36%  "
$ws.Range("D6").Value2 = "1.000"
$ws.Range("E6").Value2 = "  -0.04%  "
$ws.Range("E7").Value2 = "  -1.22%  "
$ws.Range("D8").Value2 = "0.4096"
$ws.Range("E8").Value2 = "  +2.97%  "
$ws.Range("D9").Value2 = "47.71"
$ws.Range("E9").Value2 = "  -0.11%  "
$ws.Range("D10").Value2 = "0.08014"
$ws.Range("E10").Value2 = "  -0.21%  "
$ws.Range("D11").Value2 = "1.009"
$ws.Range("E11").Value2 = "  +0.82%  "
$ws.Range("D12").Value2 = "21.92"
$ws.Range("E12").Value2 = "  -0.51%  "
$ws.Range("D13").Value2 = "1.910.09"
$ws.Range("E13").Value2 = "  +1.78%  "
$ws.Range("D14").Value2 = "5.956"
$ws.Range("E14").Value2 = "  -1.42%  "
$ws.Range("D15").Value2 = "7.109"
$ws.Range("E15").Value2 = "  -2.21%  "
$ws.Range("D16").Value2 = "89.33"
$ws.Range("E16").Value2 = "  -1.95%  "
$ws.Range("D17").Value2 = "1.001"
$ws.Range("E17").Value2 = "  +0.00%  "
$ws.Range("D18").Value2 = "0.00001033"
$ws.Range("E18").Value2 = "  -1.05%  "
$ws.Range("D19").Value2 = "0.06583"
$ws.Range("E19").Value2 = "  -0.76%  "
$ws.Range("D20").Value2 = "17.56"
$ws.Range("E20").Value2 = "  -0.08%  "
$ws.Range("D21").Value2 = "0.9999"
$ws.Range("E21").Value2 = "  -0.09%  "
$ws.Range("D22").Value2 = "29.169.56"
$ws.Range("E22").Value2 = "  +2.64%  "
$ws.Range("D23").Value2 = "5.442"
$ws.Range("E23").Value2 = "  -0.61%  "
$ws.Range("D24").Value2 = "11.31"
$ws.Range("E24").Value2 = "  +2.26%  "
$ws.Range("D25").Value2 = "2.235"
$ws.Range("E25").Value2 = "  -1.07%  "
$ws.Range("D26").Value2 = "2.132.80"
$ws.Range("E26").Value2 = "  +1.62%  "
$ws.Range("D27").Value2 = "157.95"
$ws.Range("E27").Value2 = "  -1.69%  "
$ws.Range("D28").Value2 = "19.75"
$ws.Range("E28").Value2 = "  -0.13%  "
$ws.Range("D29").Value2 = "2.118"
$ws.Range("E29").Value2 = "  -0.32%  "
$ws.Range("D30").Value2 = "5.444"
$ws.Range("E30").Value2 = "  -1.14%  "
$ws.Range("D31").Value2 = "118.61"
$ws.Range("E31").Value2 = "  -1.40%  "
$ws.Range("D32").Value2 = "0.9897"
$ws.Range("E32").Value2 = "  +1.33%  "
$ws.Range("D33").Value2 = "0.09422"
$ws.Range("E33").Value2 = "  -0.95%  "
$ws.Range("D34").Value2 = "1.436"
$ws.Range("E34").Value2 = "  +4.15%  "
$ws.Range("D35").Value2 = "3.593"
$ws.Range("E35").Value2 = "  +0.13%  "
$ws.Range("D36").Value2 = "5.324"
$ws.Range("E36").Value2 = "  -0.51%  "
$ws.Range("D37").Value2 = "0.06109"
$ws.Range("E37").Value2 = "  -0.06%  "
$ws.Range("D38").Value2 = "0.02254"
$ws.Range("E38").Value2 = "  +0.05%  "
$ws.Range("D39").Value2 = "8.388"
$ws.Range("E39").Value2 = "  +0.20%  "
$ws.Range("D40").Value2 = "1.176"
$ws.Range("E40").Value2 = "  -0.60%  "
$ws.Range("D41").Value2 = "0.5822"
$ws.Range("E41").Value2 = "  -2.22%  "
$ws.Range("D42").Value2 = "0.9996"
$ws.Range("E42").Value2 = "  -0.07%  "
$ws.Range("D43").Value2 = "10.21"
$ws.Range("E43").Value2 = "  -1.20%  "
$ws.Range("D44").Value2 = "0.1829"
$ws.Range("E44").Value2 = "  -2.80%  "
$ws.Range("D45").Value2 = "1.266"
$ws.Range("E45").Value2 = "  -1.11%  "
$ws.Range("D46").Value2 = "2.366"
$ws.Range("E46").Value2 = "  +14.95%  "
$ws.Range("D47").Value2 = "12.11"
$ws.Range("E47").Value2 = "  -0.67%  "
$ws.Range("D48").Value2 = "0.5518"
$ws.Range("E48").Value2 = "  -1.33%  "
$ws.Range("D49").Value2 = "1.925"
$ws.Range("E49").Value2 = "  -1.58%  "
$ws.Range("D50").Value2 = "0.07075"
$ws.Range("E50").Value2 = "  +1.80%  "
$ws.Range("D51").Value2 = "47.88"
$ws.Range("E51").Value2 = "  +21.93%  "
